$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteValues = -4163

function Set-TextValue($ws, $addr, $text) {
    $cell = $ws.Range($addr)
    $escaped = $text -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial($xlPasteValues)
}

Set-TextValue $ws "D2" "29.360.36"
Set-TextValue $ws "E2" "  +0.16%  "
Set-TextValue $ws "D3" "1.884.05"
Set-TextValue $ws "E3" "  +0.42%  "
Set-TextValue $ws "D4" "0.9997"
Set-TextValue $ws "E4" "  +0.00%  "
Set-TextValue $ws "D5" "0.7138"
Set-TextValue $ws "E5" "  +0.26%  "
Set-TextValue $ws "D6" "242.29"
Set-TextValue $ws "E6" "  -0.05%  "
Set-TextValue $ws "D7" "0.9999"
Set-TextValue $ws "E7" "  +0.02%  "
Set-TextValue $ws "D8" "0.08093"
Set-TextValue $ws "E8" "  +4.17%  "
Set-TextValue $ws "E9" "  +0.80%  "
Set-TextValue $ws "D10" "25.34"
Set-TextValue $ws "E10" "  +1.32%  "
Set-TextValue $ws "D11" "0.08351"
Set-TextValue $ws "E11" "  -1.31%  "
Set-TextValue $ws "D12" "1.877.54"
Set-TextValue $ws "E12" "  +0.03%  "
Set-TextValue $ws "D13" "0.7214"
Set-TextValue $ws "E13" "  +1.44%  "
Set-TextValue $ws "E14" "  +0.81%  "
Set-TextValue $ws "D15" "92.14"
Set-TextValue $ws "E15" "  +0.85%  "
Set-TextValue $ws "D16" "6.280"
Set-TextValue $ws "E16" "  +4.81%  "
Set-TextValue $ws "D17" "0.000008499"
Set-TextValue $ws "E17" "  +2.70%  "
Set-TextValue $ws "D18" "29.369.97"
Set-TextValue $ws "E18" "  +0.20%  "
Set-TextValue $ws "D19" "241.56"
Set-TextValue $ws "E19" "  -0.32%  "
Set-TextValue $ws "D20" "13.24"
Set-TextValue $ws "E20" "  +0.17%  "
Set-TextValue $ws "D21" "2.132.11"
Set-TextValue $ws "E21" "  +0.22%  "
Set-TextValue $ws "D22" "0.9996"
Set-TextValue $ws "E22" "  -0.03%  "
Set-TextValue $ws "D23" "7.810"
Set-TextValue $ws "E23" "  -0.23%  "
Set-TextValue $ws "D24" "1.0000"
Set-TextValue $ws "E24" "  +0.01%  "
Set-TextValue $ws "E25" "  -0.85%  "
Set-TextValue $ws "D26" "163.38"
Set-TextValue $ws "E26" "  +0.40%  "
Set-TextValue $ws "D27" "9.075"
Set-TextValue $ws "E27" "  +0.67%  "
Set-TextValue $ws "D28" "18.57"
Set-TextValue $ws "E28" "  +0.39%  "
Set-TextValue $ws "D29" "1.505"
Set-TextValue $ws "E29" "  -0.72%  "
Set-TextValue $ws "D30" "4.428"
Set-TextValue $ws "E30" "  +0.51%  "
Set-TextValue $ws "D31" "4.347"
Set-TextValue $ws "E31" "  +0.16%  "
Set-TextValue $ws "D32" "1.220"
Set-TextValue $ws "E32" "  -3.70%  "
Set-TextValue $ws "D33" "0.05383"
Set-TextValue $ws "E33" "  +2.41%  "
Set-TextValue $ws "E34" "  +1.82%  "
Set-TextValue $ws "D35" "1.181"
Set-TextValue $ws "E35" "  +0.69%  "
Set-TextValue $ws "D36" "0.7507"
Set-TextValue $ws "E36" "  +0.55%  "
Set-TextValue $ws "D37" "2.698"
Set-TextValue $ws "E37" "  +0.56%  "
Set-TextValue $ws "E38" "  +0.94%  "
Set-TextValue $ws "D39" "1.284.66"
Set-TextValue $ws "E39" "  +9.66%  "
Set-TextValue $ws "D40" "2.743"
Set-TextValue $ws "E40" "  +0.88%  "
Set-TextValue $ws "D41" "6.569"
Set-TextValue $ws "E41" "  +3.15%  "
Set-TextValue $ws "D42" "73.60"
Set-TextValue $ws "E42" "  +0.85%  "
Set-TextValue $ws "D43" "0.8920"
Set-TextValue $ws "E43" "  +0.60%  "
Set-TextValue $ws "D44" "110.29"
Set-TextValue $ws "D45" "1.0000"
Set-TextValue $ws "E45" "  +0.04%  "
Set-TextValue $ws "E46" "  +6.56%  "
Set-TextValue $ws "D47" "2.020.99"
Set-TextValue $ws "E47" "  +0.09%  "
Set-TextValue $ws "D48" "1.807"
Set-TextValue $ws "E48" "  -0.62%  "
Set-TextValue $ws "D49" "0.5207"
Set-TextValue $ws "E49" "  +0.19%  "
Set-TextValue $ws "D50" "9.482"
Set-TextValue $ws "E50" "  +0.99%  "
Set-TextValue $ws "D51" "0.4369"
Set-TextValue $ws "E51" "  +1.57%  "
